# Update loading_percent values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; Col="B"; Val=11.32082412494098 },
    @{ Row=2; Col="C"; Val=8.178640804395142 },
    @{ Row=2; Col="E"; Val=11.98333633642644 },
    @{ Row=2; Col="F"; Val=16.86991607391245 },
    @{ Row=2; Col="G"; Val=3.60034719021465 },
    @{ Row=2; Col="M"; Val=14.03781326338859 },
    @{ Row=2; Col="N"; Val=16.46225587529685 },
    @{ Row=2; Col="O"; Val=18.40487989437006 },
    @{ Row=3; Col="B"; Val=10.78088386508945 },
    @{ Row=3; Col="C"; Val=7.846599462231715 },
    @{ Row=3; Col="E"; Val=11.85919753509363 },
    @{ Row=3; Col="F"; Val=15.89584955866815 },
    @{ Row=3; Col="G"; Val=3.602421833494636 },
    @{ Row=3; Col="M"; Val=13.77102802358012 },
    @{ Row=3; Col="N"; Val=16.51719901889272 },
    @{ Row=3; Col="O"; Val=18.44638975706748 },
    @{ Row=4; Col="B"; Val=10.43639917922075 },
    @{ Row=4; Col="C"; Val=7.633970861245526 },
    @{ Row=4; Col="E"; Val=11.78686436687489 },
    @{ Row=4; Col="F"; Val=15.26997757108489 },
    @{ Row=4; Col="G"; Val=3.603762645725069 },
    @{ Row=4; Col="M"; Val=13.60789689591619 },
    @{ Row=4; Col="N"; Val=16.55283662722829 },
    @{ Row=4; Col="O"; Val=18.47784496985784 },
    @{ Row=5; Col="B"; Val=10.29294343273274 },
    @{ Row=5; Col="C"; Val=7.5452005243435 },
    @{ Row=5; Col="E"; Val=11.7583987950543 },
    @{ Row=5; Col="F"; Val=15.008197319934 },
    @{ Row=5; Col="G"; Val=3.604325933441762 },
    @{ Row=5; Col="M"; Val=13.54168350323451 },
    @{ Row=5; Col="N"; Val=16.56783861421201 },
    @{ Row=5; Col="O"; Val=18.49215668637765 },
    @{ Row=6; Col="B"; Val=10.26894288224622 },
    @{ Row=6; Col="C"; Val=7.530334603094911 },
    @{ Row=6; Col="E"; Val=11.75373407070418 },
    @{ Row=6; Col="F"; Val=14.96433081551589 },
    @{ Row=6; Col="G"; Val=3.604420489022538 },
    @{ Row=6; Col="M"; Val=13.53070751947653 },
    @{ Row=6; Col="N"; Val=16.57035866607288 },
    @{ Row=6; Col="O"; Val=18.49462311351134 },
    @{ Row=7; Col="B"; Val=10.43447667420977 },
    @{ Row=7; Col="C"; Val=7.632782155096564 },
    @{ Row=7; Col="E"; Val=11.78647633593462 },
    @{ Row=7; Col="F"; Val=15.26647399323133 },
    @{ Row=7; Col="G"; Val=3.60377017393694 },
    @{ Row=7; Col="M"; Val=13.60700272579879 },
    @{ Row=7; Col="N"; Val=16.55303700695334 },
    @{ Row=7; Col="O"; Val=18.47803194577422 },
    @{ Row=8; Col="B"; Val=11.13744381004496 },
    @{ Row=8; Col="C"; Val=8.066013783817871 },
    @{ Row=8; Col="E"; Val=11.93974675222471 },
    @{ Row=8; Col="F"; Val=16.5399640634477 },
    @{ Row=8; Col="G"; Val=3.601048660151961 },
    @{ Row=8; Col="M"; Val=13.94574131977208 },
    @{ Row=8; Col="N"; Val=16.48080606360775 },
    @{ Row=8; Col="O"; Val=18.41794967669124 },
    @{ Row=9; Col="B"; Val=12.4064251930495 },
    @{ Row=9; Col="C"; Val=8.84318426732953 },
    @{ Row=9; Col="E"; Val=12.26959709458786 },
    @{ Row=9; Col="F"; Val=19.00274580682531 },
    @{ Row=9; Col="G"; Val=3.596240668512143 },
    @{ Row=9; Col="M"; Val=14.61132434598412 },
    @{ Row=9; Col="N"; Val=16.35420900324438 },
    @{ Row=9; Col="O"; Val=18.34777648412739 },
    @{ Row=10; Col="B"; Val=13.26461281928683 },
    @{ Row=10; Col="C"; Val=9.36675546383103 },
    @{ Row=10; Col="E"; Val=12.52755424217144 },
    @{ Row=10; Col="F"; Val=20.67494806633232 },
    @{ Row=10; Col="G"; Val=3.593027143152062 },
    @{ Row=10; Col="M"; Val=15.09587707362791 },
    @{ Row=10; Col="N"; Val=16.27030751230372 },
    @{ Row=10; Col="O"; Val=18.32562644954232 },
    @{ Row=11; Col="B"; Val=13.63776879413103 },
    @{ Row=11; Col="C"; Val=9.594113689800357 },
    @{ Row=11; Col="E"; Val=12.64778996309311 },
    @{ Row=11; Col="F"; Val=21.3917225636224 },
    @{ Row=11; Col="G"; Val=3.591633725427643 },
    @{ Row=11; Col="M"; Val=15.31429657209859 },
    @{ Row=11; Col="N"; Val=16.23410331943791 },
    @{ Row=11; Col="O"; Val=18.32199628996652 },
    @{ Row=12; Col="B"; Val=13.77651913623583 },
    @{ Row=12; Col="C"; Val=9.678617313691101 },
    @{ Row=12; Col="E"; Val=12.6936919476929 },
    @{ Row=12; Col="F"; Val=21.65686569030329 },
    @{ Row=12; Col="G"; Val=3.591115857837734 },
    @{ Row=12; Col="M"; Val=15.39663571509501 },
    @{ Row=12; Col="N"; Val=16.22067506530021 },
    @{ Row=12; Col="O"; Val=18.3215522604234 },
    @{ Row=13; Col="B"; Val=13.74675151696102 },
    @{ Row=13; Col="C"; Val=9.660489271145886 },
    @{ Row=13; Col="E"; Val=12.6837903269487 },
    @{ Row=13; Col="F"; Val=21.60004134736742 },
    @{ Row=13; Col="G"; Val=3.591226955305846 },
    @{ Row=13; Col="M"; Val=15.37892030964901 },
    @{ Row=13; Col="N"; Val=16.22355457315284 },
    @{ Row=13; Col="O"; Val=18.32160645006493 },
    @{ Row=14; Col="B"; Val=13.6492354817687 },
    @{ Row=14; Col="C"; Val=9.601097948470851 },
    @{ Row=14; Col="E"; Val=12.6515591596143 },
    @{ Row=14; Col="F"; Val=21.4136618050453 },
    @{ Row=14; Col="G"; Val=3.591590924244036 },
    @{ Row=14; Col="M"; Val=15.32107857027153 },
    @{ Row=14; Col="N"; Val=16.23299293188917 },
    @{ Row=14; Col="O"; Val=18.3219410899299 },
    @{ Row=15; Col="B"; Val=13.5891692107447 },
    @{ Row=15; Col="C"; Val=9.564510738284532 },
    @{ Row=15; Col="E"; Val=12.63186367789013 },
    @{ Row=15; Col="F"; Val=21.29868154950795 },
    @{ Row=15; Col="G"; Val=3.591815139102243 },
    @{ Row=15; Col="M"; Val=15.28559805280501 },
    @{ Row=15; Col="N"; Val=16.23881083415005 },
    @{ Row=15; Col="O"; Val=18.32226735781237 },
    @{ Row=16; Col="B"; Val=13.23987243692863 },
    @{ Row=16; Col="C"; Val=9.351676145368865 },
    @{ Row=16; Col="E"; Val=12.51975093427618 },
    @{ Row=16; Col="F"; Val=20.62722412089977 },
    @{ Row=16; Col="G"; Val=3.593119578930047 },
    @{ Row=16; Col="M"; Val=15.08155527415119 },
    @{ Row=16; Col="N"; Val=16.27271297409904 },
    @{ Row=16; Col="O"; Val=18.32599373526277 },
    @{ Row=17; Col="B"; Val=13.02111596062233 },
    @{ Row=17; Col="C"; Val=9.218310628581449 },
    @{ Row=17; Col="E"; Val=12.45168222033107 },
    @{ Row=17; Col="F"; Val=20.20408069597325 },
    @{ Row=17; Col="G"; Val=3.593937301062434 },
    @{ Row=17; Col="M"; Val=14.95580732997229 },
    @{ Row=17; Col="N"; Val=16.29401304292017 },
    @{ Row=17; Col="O"; Val=18.32993362890987 },
    @{ Row=18; Col="B"; Val=12.89367415845148 },
    @{ Row=18; Col="C"; Val=9.140585948206216 },
    @{ Row=18; Col="E"; Val=12.4128055773666 },
    @{ Row=18; Col="F"; Val=19.95656407809801 },
    @{ Row=18; Col="G"; Val=3.594414077277381 },
    @{ Row=18; Col="M"; Val=14.88329632953616 },
    @{ Row=18; Col="N"; Val=16.3064490930903 },
    @{ Row=18; Col="O"; Val=18.33280631117665 },
    @{ Row=19; Col="B"; Val=12.85024908445067 },
    @{ Row=19; Col="C"; Val=9.114096357069581 },
    @{ Row=19; Col="E"; Val=12.39969113033285 },
    @{ Row=19; Col="F"; Val=19.87204792380568 },
    @{ Row=19; Col="G"; Val=3.594576613842788 },
    @{ Row=19; Col="M"; Val=14.85871633458606 },
    @{ Row=19; Col="N"; Val=16.31069149253989 },
    @{ Row=19; Col="O"; Val=18.33388300245791 },
    @{ Row=20; Col="B"; Val=13.04457110783554 },
    @{ Row=20; Col="C"; Val=9.232613088631659 },
    @{ Row=20; Col="E"; Val=12.45890014160845 },
    @{ Row=20; Col="F"; Val=20.24955283636154 },
    @{ Row=20; Col="G"; Val=3.593849586577003 },
    @{ Row=20; Col="M"; Val=14.96921308403487 },
    @{ Row=20; Col="N"; Val=16.29172649248232 },
    @{ Row=20; Col="O"; Val=18.32945141427484 },
    @{ Row=21; Col="B"; Val=13.67794821332661 },
    @{ Row=21; Col="C"; Val=9.618586088395833 },
    @{ Row=21; Col="E"; Val=12.66101651124716 },
    @{ Row=21; Col="F"; Val=21.46857628470577 },
    @{ Row=21; Col="G"; Val=3.591483752464819 },
    @{ Row=21; Col="M"; Val=15.33807881619828 },
    @{ Row=21; Col="N"; Val=16.23021302381057 },
    @{ Row=21; Col="O"; Val=18.32181751567877 },
    @{ Row=22; Col="B"; Val=14.07697343977191 },
    @{ Row=22; Col="C"; Val=9.861549367040869 },
    @{ Row=22; Col="E"; Val=12.795254924293 },
    @{ Row=22; Col="F"; Val=22.22866616901552 },
    @{ Row=22; Col="G"; Val=3.589994580239181 },
    @{ Row=22; Col="M"; Val=15.57695267360308 },
    @{ Row=22; Col="N"; Val=16.19165079502401 },
    @{ Row=22; Col="O"; Val=18.32225387420913 },
    @{ Row=23; Col="B"; Val=13.8653936805895 },
    @{ Row=23; Col="C"; Val=9.732736144206058 },
    @{ Row=23; Col="E"; Val=12.72342789503219 },
    @{ Row=23; Col="F"; Val=21.82633154458858 },
    @{ Row=23; Col="G"; Val=3.590784177286769 },
    @{ Row=23; Col="M"; Val=15.44968891210884 },
    @{ Row=23; Col="N"; Val=16.21208234207482 },
    @{ Row=23; Col="O"; Val=18.32152354081267 },
    @{ Row=24; Col="B"; Val=13.03397224951065 },
    @{ Row=24; Col="C"; Val=9.226150217520113 },
    @{ Row=24; Col="E"; Val=12.45563611562148 },
    @{ Row=24; Col="F"; Val=20.22900810905287 },
    @{ Row=24; Col="G"; Val=3.593889221556011 },
    @{ Row=24; Col="M"; Val=14.96315301162099 },
    @{ Row=24; Col="N"; Val=16.29275964889727 },
    @{ Row=24; Col="O"; Val=18.32966753116596 },
    @{ Row=25; Col="B"; Val=12.07569938010259 },
    @{ Row=25; Col="C"; Val=8.641061052863572 },
    @{ Row=25; Col="E"; Val=12.17745305949718 },
    @{ Row=25; Col="F"; Val=18.34778573295695 },
    @{ Row=25; Col="G"; Val=3.597485101547926 },
    @{ Row=25; Col="M"; Val=14.43167945684994 },
    @{ Row=25; Col="N"; Val=16.38685263184929 },
    @{ Row=25; Col="O"; Val=18.36161918829445 }
)

foreach ($item in $updates) {
    $addr = $item.Col + $item.Row
    $ws.Range($addr).Value = $item.Val
}

